$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.933.48'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.60%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.136.77'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.53%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '612.14'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.16%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.31'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.90%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.133.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.67%  '

$ws.Range("E9").Value = '  -4.03%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.30'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -8.35%  '

$ws.Range("E12").Value = '  -5.56%  '

$ws.Range("E13").Value = '  -7.54%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.24'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -9.80%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.653.17'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.41%  '

$ws.Range("E16").Value = '  +1.17%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.964.22'
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.137.56'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.27%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.84'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -8.22%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '475.51'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.79%  '

$ws.Range("E21").Value = '  -5.29%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.703'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -7.04%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.75'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.14%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.56'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -7.68%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.35'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.62%  '

$ws.Range("E26").Value = '  -0.08%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.79'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -7.69%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.40'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -8.40%  '

$ws.Range("E29").Value = '  -9.71%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.73'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.28%  '

$ws.Range("E31").Value = '  -10.07%  '

$ws.Range("B32").Value = 'Stacks'
$ws.Range("C32").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.73'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.29%  '

$ws.Range("B33").Value = 'FirstDigitalUSD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.05%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.14'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.48%  '

$ws.Range("E35").Value = '  -2.58%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.94'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -8.12%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '53.20'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.25%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0731'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.34%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '459.58'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.21%  '

$ws.Range("E40").Value = '  -13.31%  '

$ws.Range("E41").Value = '  -7.09%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.120'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.20%  '

$ws.Range("E43").Value = '  -5.74%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.837.05'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.14%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.265'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -9.74%  '

$ws.Range("E46").Value = '  -11.20%  '

$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '26.28'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -9.28%  '

$ws.Range("B49").Value = 'ThetaToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.36'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.57%  '

$ws.Range("E50").Value = '  -5.12%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '118.28'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.02%  '
